# Apply updated crypto price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/Volume columns hold plain text (e.g. "52.565.55", "  +0.63%  ").
# Force text storage first so numeric-looking values are not coerced to
# real numbers, then drop back to the default style once written.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '52.565.55'
$ws.Range("E2").Value = '  +0.63%  '

$ws.Range("D3").Value = '2.982.75'
$ws.Range("E3").Value = '  +2.77%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").Value = '358.25'
$ws.Range("E5").Value = '  +1.29%  '

$ws.Range("D6").Value = '110.09'
$ws.Range("E6").Value = '  -3.27%  '

$ws.Range("D7").Value = '0.574'
$ws.Range("E7").Value = '  +3.07%  '

$ws.Range("D8").Value = '0.997'
$ws.Range("E8").Value = '  -0.33%  '

$ws.Range("D9").Value = '0.632'
$ws.Range("E9").Value = '  +1.02%  '

$ws.Range("D10").Value = '39.23'
$ws.Range("E10").Value = '  -1.98%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.0879'
$ws.Range("E11").Value = '  +1.76%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.138'
$ws.Range("E12").Value = '  +1.62%  '

$ws.Range("D13").Value = '19.56'
$ws.Range("E13").Value = '  -1.45%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.430.98'
$ws.Range("E14").Value = '  +2.05%  '

$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '7.82'
$ws.Range("E15").Value = '  +0.87%  '

$ws.Range("D16").Value = '2.972.78'
$ws.Range("E16").Value = '  +1.96%  '

$ws.Range("D17").Value = '0.988'
$ws.Range("E17").Value = '  -0.93%  '

$ws.Range("D18").Value = '52.428.61'
$ws.Range("E18").Value = '  +0.31%  '

$ws.Range("D19").Value = '3.49'
$ws.Range("E19").Value = '  +5.83%  '

$ws.Range("D20").Value = '7.69'
$ws.Range("E20").Value = '  +0.55%  '

$ws.Range("D21").Value = '14.01'
$ws.Range("E21").Value = '  -1.09%  '

$ws.Range("D22").Value = '0.0₃0987'
$ws.Range("E22").Value = '  +0.93%  '

$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value = '272.51'
$ws.Range("E23").Value = '  +1.04%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '70.61'
$ws.Range("E24").Value = '  -0.33%  '

$ws.Range("D25").Value = '2.81'
$ws.Range("E25").Value = '  +0.78%  '

$ws.Range("D26").Value = '0.180'
$ws.Range("E26").Value = '  +3.14%  '

$ws.Range("D27").Value = '7.91'
$ws.Range("E27").Value = '  +19.81%  '

$ws.Range("D28").Value = '27.34'
$ws.Range("E28").Value = '  +2.12%  '

$ws.Range("E29").Value = '  +0.16%  '

$ws.Range("D30").Value = '0.109'
$ws.Range("E30").Value = '  +4.41%  '

$ws.Range("D31").Value = '10.53'
$ws.Range("E31").Value = '  -0.76%  '

$ws.Range("D32").Value = '38.03'
$ws.Range("E32").Value = '  +1.10%  '

$ws.Range("D33").Value = '6.16'
$ws.Range("E33").Value = '  -1.22%  '

$ws.Range("E34").Value = '  +10.85%  '

$ws.Range("D35").Value = '52.47'
$ws.Range("E35").Value = '  -1.34%  '

$ws.Range("D36").Value = '0.0444'
$ws.Range("E36").Value = '  -1.12%  '

$ws.Range("E37").Value = '  -0.08%  '

$ws.Range("D38").Value = '3.25'
$ws.Range("E38").Value = '  -2.36%  '

$ws.Range("D39").Value = '2.03'
$ws.Range("E39").Value = '  -1.28%  '

$ws.Range("D40").Value = '18.28'
$ws.Range("E40").Value = '  -3.34%  '

$ws.Range("D41").Value = '2.75'
$ws.Range("E41").Value = '  -0.44%  '

$ws.Range("E42").Value = '  +2.91%  '

$ws.Range("D43").Value = '23.88'
$ws.Range("E43").Value = '  +4.00%  '

$ws.Range("D44").Value = '118.86'
$ws.Range("E44").Value = '  -0.87%  '

$ws.Range("E45").Value = '  -0.81%  '

$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = '3.47'
$ws.Range("E46").Value = '  -1.31%  '

$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '2.47'
$ws.Range("E47").Value = '  -5.67%  '

$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '2.149.95'
$ws.Range("E48").Value = '  -1.36%  '

$ws.Range("B49").Value = 'BEAM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range("D49").Value = '0.0357'
$ws.Range("E49").Value = '  +1.63%  '

$ws.Range("B50").Value = 'TheGraph'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D50").Value = '0.248'
$ws.Range("E50").Value = '  -5.42%  '

$ws.Range("D51").Value = '0.920'
$ws.Range("E51").Value = '  -3.64%  '

# Restore the default (unstyled) look for the cells we just wrote.
$ws.Range("D2:E51").Style = "Normal"
